# Edit script: adds "Burmester" and "Burmester-modernized" columns to the
# Ps 118-10 comparison table (after "Edited", before "AI"), resizes all
# columns, and fills in the Burmester translation text for each verse.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- 1. Insert the two new columns before the current column 3 ("AI") ---
$t.Columns.Add($t.Columns.Item(3)) | Out-Null
$t.Columns.Add($t.Columns.Item(3)) | Out-Null

# --- 2. Burmester column text per row (row 1 = header) ---
$burmester = @{
    1 = "Burmester"
    2 = "Thy hands have formed me and created me: teach me and I shall know Thy commandments."
    3 = "Those who fear Thee shall see me, and they shall be glad, for I have hoped in Thy word."
    4 = "75 I have Known, Lord, that Thy judgments are righteous: in truth Thou hast humbled me."
    5 = "76 Let Thy mercy come upon me to comfort me, and Thy word to Thy servant."
    6 = "Let Thy tender-mercies come to me, and I shall live: for Thy Law is my meditation:"
    7 = "Let the haughty be ashamed, for unjustly have they transgressed against me: but I, I shall be continuing in Thy commandments."
    8 = "Let those who fear Thee and those who know Thy wonders turn unto me."
    9 = "Let my heart be pure in Thy truth, that I be not ashamed."
}

# --- 3. Burmester-modernized column text per row (only header has text) ---
$burmesterModernized = @{
    1 = "Burmester-modernized"
}

$rowCount = $t.Rows.Count
for ($r = 1; $r -le $rowCount; $r++) {
    $cellB = $t.Cell($r, 3)
    $cellB.Range.Text = $burmester[$r]

    if ($burmesterModernized.ContainsKey($r)) {
        $cellBM = $t.Cell($r, 4)
        $cellBM.Range.Text = $burmesterModernized[$r]
    }
}

# --- 4. Resize every cell in every row to the target percentage widths ---
$pctWidths = @(532, 518, 388, 388, 529, 519, 525, 534, 534, 533)
$colCount = $t.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.PreferredWidthType = 2
        $cell.PreferredWidth = $pctWidths[$c - 1] / 20.0
    }
}

Write-Host "Done"
